# Daily attendance processing - 2025-11-11 13:54:37
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# holds a comma-separated list of recorder identities per attendance
# session (e.g. "System, dnasr281@gmail.com"). This pass reverses the
# ordering of the entries in each such list (cells with only a single
# entry are left unchanged, since reversing a one-item list is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 157
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $txt = $cell.Text

    if ($txt -eq $null -or $txt -eq "") {
        continue
    }

    $parts = $txt -split ", "

    if ($parts.Length -gt 1) {
        $reversed = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = $reversed -join ", "
    }
}

Write-Output "Recorded By column reversed through row $lastRow"
